$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "last updated" timestamp banner ---
$ws.Range("A1").Value2 = "Datos actualizados a 16 de Abril de 2020 a las 05:52"

# --- Update India's stats (row 23) ---
$ws.Range("B23").Value2 = 12380
$ws.Range("C23").Value2 = 10
$ws.Range("E23").Value2 = 10450

# --- Re-insert "Birmania" into its alphabetically-correct spot ---
# Currently (row 140..147): Etiopia, Bermudas, Togo, Gabon, Somalia, Liechtenstein,
# Barbados, Birmania. Birmania needs to move up to row 140 (just before Etiopia),
# pushing the other seven rows down by one.
$ws.Range("A140:H140").Insert(1) | Out-Null
$ws.Range("A148:H148").Cut($ws.Range("A140:H140")) | Out-Null
$ws.Range("A148:H148").Delete(1) | Out-Null

# --- Refresh Birmania's case numbers for row 140 ---
$ws.Range("B140").Value2 = 85
$ws.Range("C140").Value2 = 11
$ws.Range("D140").Value2 = 2
$ws.Range("E140").Value2 = 79
$ws.Range("F140").Value2 = 0
$ws.Range("G140").Value2 = 0
$ws.Range("H140").Value2 = 4
